# Update contents of production machine
# - Rename "Alumno" header to "Alumnado"
# - Replace old "imwpto.me" domain value with "vps.claseando.es"
# - Add a new "Host" row (claseando / cloud) styled like the existing rows
#   but with a new accent (green) fill color

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Header C1: "Alumno" -> "Alumnado"
$ws.Range("C1").Value = "Alumnado"

# 2. Domain value B3: "imwpto.me" -> "vps.claseando.es"
$ws.Range("B3").Value = "vps.claseando.es"

# 3. New row 4 content: Host / claseando / cloud
$ws.Range("A4").Value = "Host"
$ws.Range("B4").Value = "claseando"
$ws.Range("C4").Value = "cloud"

# Match formatting of the row above (font/border) for the new row
$ws.Range("A3").Copy()
$ws.Range("A4").PasteSpecial(-4122)
$ws.Range("B3:C3").Copy()
$ws.Range("B4:C4").PasteSpecial(-4122)

# Give the new data row its own accent fill color (Accent6 / green)
$ws.Range("B4:C4").Interior.ThemeColor = 10

# Widen column B to fit the longer "vps.claseando.es" text
$ws.Columns.Item(2).ColumnWidth = 16.33

# Restore active selection
$ws.Range("B7").Select()

Write-Output "done"
